$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string that must be preserved
# verbatim (e.g. trailing zeros like "1.00"). Force Text format first so
# Excel does not silently reinterpret/renormalize them as numbers.
$textCells = @('D9', 'D30', 'D32', 'D39', 'D42', 'D43', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '91.845.20'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '3.122.37'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '243.47'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('D7').Value = '1.11'
$ws.Range('E7').Value = '  -3.64%  '
$ws.Range('E8').Value = '  +3.28%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '3.122.61'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = '0.761'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '0.205'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '0.0000253'
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('D14').Value = '35.55'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '5.63'
$ws.Range('E15').Value = '  +1.73%  '
$ws.Range('D16').Value = '91.644.79'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '3.188.48'
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').Value = '3.81'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').Value = '14.98'
$ws.Range('E20').Value = '  +1.86%  '
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').Value = '458.05'
$ws.Range('E22').Value = '  +1.12%  '
$ws.Range('E23').Value = '  -6.66%  '
$ws.Range('D24').Value = '9.19'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '5.98'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').Value = '89.58'
$ws.Range('E26').Value = '  -3.81%  '
$ws.Range('D27').Value = '11.78'
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('D28').Value = '1.39'
$ws.Range('E28').Value = '  +37.81%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.141'
$ws.Range('E31').Value = '  +13.30%  '
$ws.Range('D32').Value = '0.230'
$ws.Range('E32').Value = '  -1.05%  '
$ws.Range('E33').Value = '  -5.87%  '
$ws.Range('D34').Value = '9.41'
$ws.Range('E34').Value = '  +2.45%  '
$ws.Range('D35').Value = '0.175'
$ws.Range('E35').Value = '  +7.68%  '
$ws.Range('B36').Value = 'PancakeSwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D36').Value = '2.23'
$ws.Range('E36').Value = '  +14.99%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '26.51'
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = '7.53'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = '492.40'
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').Value = '3.84'
$ws.Range('E41').Value = '  -8.50%  '
$ws.Range('D42').Value = '0.440'
$ws.Range('D43').Value = '3.40'
$ws.Range('E43').Value = '  -6.85%  '
$ws.Range('D44').Value = '22.19'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D46').Value = '0.709'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').Value = '156.22'
$ws.Range('E48').Value = '  -0.86%  '
$ws.Range('D49').Value = '1.36'
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('D50').Value = '4.50'
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('D51').Value = '0.0329'
$ws.Range('E51').Value = '  +1.42%  '

Write-Host "Applied updates to cryptos sheet"
